$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "Price" (column D) values are numeric-looking strings (e.g. "67.60",
# "0.809") that must stay exact text, matching the source inlineStr cells --
# not be auto-coerced by Excel into real numbers (which would silently drop
# trailing zeros / change the stored type). Force those specific cells to
# Text format before assigning so the literal string is preserved.

$ws.Range('D2').Value = '27.205.77'
$ws.Range('E2').Value = '  +1.32%  '
$ws.Range('D3').Value = '1.650.24'
$ws.Range('E3').Value = '  +0.34%  '
$ws.Range('E4').Value = '  -0.29%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '217.37'
$ws.Range('E5').Value = '  +0.22%  '
$ws.Range('E6').Value = '  +2.61%  '
$ws.Range('E7').Value = '  -0.29%  '
$ws.Range('E8').Value = '  +1.57%  '
$ws.Range('E9').Value = '  +1.45%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.03'
$ws.Range('E10').Value = '  +1.67%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0849'
$ws.Range('E11').Value = '  +0.51%  '
$ws.Range('D12').Value = '1.882.48'
$ws.Range('E12').Value = '  +0.33%  '
$ws.Range('D13').Value = '1.654.23'
$ws.Range('E13').Value = '  +0.57%  '
$ws.Range('E14').Value = '  +0.57%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.542'
$ws.Range('E15').Value = '  +2.90%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '67.60'
$ws.Range('E16').Value = '  +2.00%  '
$ws.Range('D17').Value = '27.229.64'
$ws.Range('E17').Value = '  +1.22%  '
$ws.Range('D18').Value = '0.0₃0739'
$ws.Range('E18').Value = '  +1.27%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '219.45'
$ws.Range('E19').Value = '  +0.43%  '
$ws.Range('E20').Value = '  -0.18%  '
$ws.Range('E21').Value = '  +3.37%  '
$ws.Range('E22').Value = '  +6.35%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.43'
$ws.Range('E23').Value = '  +1.17%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.23'
$ws.Range('E24').Value = '  +0.82%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '147.69'
$ws.Range('E25').Value = '  +1.31%  '
$ws.Range('E26').Value = '  +2.82%  '
$ws.Range('E27').Value = '  -0.33%  '
$ws.Range('E28').Value = '  +0.33%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.80'
$ws.Range('E29').Value = '  -0.33%  '
$ws.Range('E30').Value = '  -0.38%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.18'
$ws.Range('E31').Value = '  -0.14%  '
$ws.Range('E32').Value = '  +0.94%  '
$ws.Range('E33').Value = '  +2.16%  '
$ws.Range('E34').Value = '  +1.77%  '
$ws.Range('D35').Value = '1.271.40'
$ws.Range('E35').Value = '  +2.08%  '
$ws.Range('E36').Value = '  +0.12%  '
$ws.Range('E37').Value = '  +1.80%  '
$ws.Range('E38').Value = '  +3.03%  '
$ws.Range('E39').Value = '  +2.49%  '
$ws.Range('E40').Value = '  -0.23%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.809'
$ws.Range('E41').Value = '  +0.15%  '
$ws.Range('E42').Value = '  +1.77%  '
$ws.Range('E43').Value = '  +6.01%  '
$ws.Range('D44').Value = '1.792.08'
$ws.Range('E44').Value = '  +0.16%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '62.15'
$ws.Range('E45').Value = '  +2.06%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '91.72'
$ws.Range('E46').Value = '  +0.33%  '
$ws.Range('E47').Value = '  +0.86%  '
$ws.Range('E48').Value = '  -0.16%  '
$ws.Range('E49').Value = '  -0.10%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.69'
$ws.Range('E50').Value = '  +1.91%  '
$ws.Range('E51').Value = '  +0.24%  '
